$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$excel.ActiveWindow.TopLeftCell = $ws2.Range("A12")
Write-Host "Window TopLeftCell addr:" $excel.ActiveWindow.TopLeftCell.Address()
